# Daily GitHub Actions refresh of the crypto price table: re-pull the
# latest price/volume figures for the existing coins, and replace three
# coins that dropped out of the ranking window (One, TigerCash, HotbitToken)
# with three that entered it (TigerCash, HotbitToken, BitKan, NitroEx, LEO,
# BTSEToken, One), shifting rows 18-24 by one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores its numbers as plain text (e.g. "2.130", "0.07420")
# so that trailing zeros / exact precision survive. A bare .Value assignment
# of a numeric-looking string gets auto-coerced to a real number by Excel
# (dropping trailing zeros, switching tiny values to scientific notation),
# so force those cells to Text format before writing the new figures.
$textCells = @(
    "D2", "D3", "D4", "D6", "D8", "D9",
    "D11", "D12", "D13", "D14", "D15", "D16", "D17",
    "D18", "D19", "D20", "D21", "D22", "D23", "D24",
    "D26",
    "D40", "D41", "D42", "D43", "D44", "D45",
    "D48"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Column D price refresh (rows whose Coin/Link/Volume stayed the same) ---
$ws.Range("D2").Value = "244.05"
$ws.Range("D3").Value = "23.27"
$ws.Range("D4").Value = "5.412"
$ws.Range("D6").Value = "3.464"
$ws.Range("D8").Value = "0.8137"
$ws.Range("D9").Value = "0.9185"
$ws.Range("D11").Value = "0.07420"
$ws.Range("D12").Value = "0.03248"
$ws.Range("D13").Value = "0.03086"
$ws.Range("D14").Value = "0.09351"
$ws.Range("D15").Value = "3.846"
$ws.Range("D16").Value = "0.001559"
$ws.Range("D17").Value = "0.04678"

# --- Rows 18-24: coin list shifted by one (new entrant pushes the rest down) ---
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "0.006059"
$ws.Range("E18").Value = "17TigerCashTCH"

$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "0.005010"
$ws.Range("E19").Value = "18HotbitTokenHTB"

$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "0.0009852"
$ws.Range("E20").Value = "19BitKanKAN"

$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "0.00007799"
$ws.Range("E21").Value = "20NitroExNTX"

$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "3.612"
$ws.Range("E22").Value = "21LEOLEO"

$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").Value = "2.130"
$ws.Range("E23").Value = "22BTSETokenBTSE"

$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").Value = "0.01129"
$ws.Range("E24").Value = "23OneONEBestin24h"

# --- Remaining column D price refreshes further down the sheet ---
$ws.Range("D26").Value = "0.1296"
$ws.Range("D40").Value = "0.03925"
$ws.Range("D41").Value = "0.006223"
$ws.Range("D42").Value = "0.1073"
$ws.Range("D43").Value = "0.002700"
$ws.Range("D44").Value = "0.007097"
$ws.Range("D45").Value = "0.00005233"
$ws.Range("D48").Value = "0.8499"
